$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.480.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07890"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.882"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.471.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.340"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.043.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.065"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9497"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09316"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.249"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.330"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.031"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5768"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5458"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.874"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06612"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("E51").Value = "  -1.38%  "
